$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update batch_size (J) and negatives (K) columns for rows 2-6
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 10).Value = 8
    $ws.Cells.Item($row, 11).Value = 4
}

# Update the active selection to match the diff (K11)
$ws.Range("K11").Select()
